$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Sudo masina / DW50K / ne pere kako treba)
$ws.Range("B2").Value = "Sudo mašina"
$ws.Range("D2").Value = "DW50K"
$ws.Range("F2").Value = "ne pere kako treba"

# Row 3 updates (Ves masina / WM-5000 / ne centrifugira)
$ws.Range("B3").Value = "Veš mašina"
$ws.Range("D3").Value = "WM-5000"
$ws.Range("F3").Value = "ne centrifugira"

# Row 4 updates (Kombinovana ves masina / WKD-300 / ne susi)
$ws.Range("B4").Value = "Kombinovana veš mašina"
$ws.Range("D4").Value = "WKD-300"
$ws.Range("F4").Value = "ne suši"

# Row 5 updates (Ugradna sudo masina / WDI-60 / ne radi)
$ws.Range("B5").Value = "Ugradna sudo mašina"
$ws.Range("D5").Value = "WDI-60"
$ws.Range("F5").Value = "ne radi"

# Row 6 updates (Frizider capitalized / ne hladi)
$ws.Range("B6").Value = "Frižider"
$ws.Range("F6").Value = "ne hladi"

# New row 7 - Marko Petrović / Šporet / Gorenje / G-500 / 11111 / ne radi ploca
$ws.Range("A7").Value = "Marko Petrović"
$ws.Range("B7").Value = "Šporet"
$ws.Range("C7").Value = "Gorenje"
$ws.Range("D7").Value = "G-500"
$ws.Range("E7").Value = "'11111"
$ws.Range("F7").Value = "ne radi ploca"
